$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2025-04-28 06:28:06"
$ws.Range("B4").Value = 202
